$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024607700635785
$ws.Range("D2").Value = 1.026149859510542
$ws.Range("E2").Value = 1.025073165548221
$ws.Range("F2").Value = 1.035154460294946
$ws.Range("I2").Value = 1.032172235627278
$ws.Range("J2").Value = 1.029781365592948
$ws.Range("K2").Value = 1.028973348598866
$ws.Range("L2").Value = 1.027899805325654
$ws.Range("M2").Value = 1.03795187841028
$ws.Range("N2").Value = 1.031243773249667

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.025576526147129
$ws.Range("D3").Value = 1.026983881884379
$ws.Range("E3").Value = 1.02589534093246
$ws.Range("F3").Value = 1.036344313573487
$ws.Range("I3").Value = 1.032353193227729
$ws.Range("J3").Value = 1.030389070285153
$ws.Range("K3").Value = 1.029614752364802
$ws.Range("L3").Value = 1.028529165815819
$ws.Range("M3").Value = 1.038950055937122
$ws.Range("N3").Value = 1.031852340952236

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.02620334923417
$ws.Range("D4").Value = 1.027523791124847
$ws.Range("E4").Value = 1.026427681397812
$ws.Range("F4").Value = 1.037114056540733
$ws.Range("I4").Value = 1.032468239768286
$ws.Range("J4").Value = 1.030781642468573
$ws.Range("K4").Value = 1.030029381982526
$ws.Range("L4").Value = 1.028936096998637
$ws.Range("M4").Value = 1.039595207164697
$ws.Range("N4").Value = 1.032245470633185

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.026466848117759
$ws.Range("D5").Value = 1.027750825824427
$ws.Range("E5").Value = 1.026651557372716
$ws.Range("F5").Value = 1.037437615841038
$ws.Range("I5").Value = 1.032516115256485
$ws.Range("J5").Value = 1.030946522838513
$ws.Range("K5").Value = 1.030203595534217
$ws.Range("L5").Value = 1.029107096704974
$ws.Range("M5").Value = 1.039866252662251
$ws.Range("N5").Value = 1.032410585152163

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.02651108970076
$ws.Range("D6").Value = 1.027788949258941
$ws.Range("E6").Value = 1.026689151802129
$ws.Range("F6").Value = 1.037491940501051
$ws.Range("I6").Value = 1.032524125009443
$ws.Range("J6").Value = 1.030974197756777
$ws.Range("K6").Value = 1.03023284106519
$ws.Range("L6").Value = 1.029135803940627
$ws.Range("M6").Value = 1.039911752074648
$ws.Range("N6").Value = 1.032438299371985

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026206870188527
$ws.Range("D7").Value = 1.027526824551283
$ws.Range("E7").Value = 1.026430672526147
$ws.Range("F7").Value = 1.037118380115092
$ws.Range("I7").Value = 1.032468881409176
$ws.Range("J7").Value = 1.030783846225663
$ws.Range("K7").Value = 1.03003171021344
$ws.Range("L7").Value = 1.0289383821973
$ws.Range("M7").Value = 1.039598829579529
$ws.Range("N7").Value = 1.032247677519863

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.024935134410239
$ws.Range("D8").Value = 1.026431670874395
$ws.Range("E8").Value = 1.025350953054901
$ws.Range("F8").Value = 1.035556612601979
$ws.Range("I8").Value = 1.032233814272464
$ws.Range("J8").Value = 1.029986877196788
$ws.Range("K8").Value = 1.029190196770295
$ws.Range("L8").Value = 1.028112563866246
$ws.Range("M8").Value = 1.038289369781221
$ws.Range("N8").Value = 1.031449576703554

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02269363840731
$ws.Range("D9").Value = 1.024503752382478
$ws.Range("E9").Value = 1.023450976384379
$ws.Range("F9").Value = 1.032803242636772
$ws.Range("I9").Value = 1.031803952933132
$ws.Range("J9").Value = 1.028577535593116
$ws.Range("K9").Value = 1.027704294184937
$ws.Range("L9").Value = 1.026655041559764
$ws.Range("M9").Value = 1.035976284034158
$ws.Range("N9").Value = 1.030038233673151

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021198964917186
$ws.Range("D10").Value = 1.023219785937863
$ws.Range("E10").Value = 1.022186141038229
$ws.Range("F10").Value = 1.030966730617533
$ws.Range("I10").Value = 1.031506887069839
$ws.Range("J10").Value = 1.027634656063815
$ws.Range("K10").Value = 1.026711676897278
$ws.Range("L10").Value = 1.025681831986985
$ws.Range("M10").Value = 1.034430409388537
$ws.Range("N10").Value = 1.029094015146768

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020551676063724
$ws.Range("D11").Value = 1.022664134867342
$ws.Range("E11").Value = 1.021638893868459
$ws.Range("F11").Value = 1.03017127083266
$ws.Range("I11").Value = 1.031375770779068
$ws.Range("J11").Value = 1.027225597394939
$ws.Range("K11").Value = 1.026281391122717
$ws.Range("L11").Value = 1.025260065873007
$ws.Range("M11").Value = 1.033760117813289
$ws.Range("N11").Value = 1.028684375567663

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020311231430183
$ws.Range("D12").Value = 1.022457789299756
$ws.Range("E12").Value = 1.021435687965769
$ws.Range("F12").Value = 1.029875765008598
$ws.Range("I12").Value = 1.031326695366954
$ws.Range("J12").Value = 1.027073537139403
$ws.Range("K12").Value = 1.026121492775224
$ws.Range("L12").Value = 1.025103349560377
$ws.Range("M12").Value = 1.033511003212293
$ws.Range("N12").Value = 1.028532099369122

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020362808195032
$ws.Range("D13").Value = 1.0225020489425
$ws.Range("E13").Value = 1.02147927332065
$ws.Range("F13").Value = 1.029939153683501
$ws.Range("I13").Value = 1.031337239082497
$ws.Range("J13").Value = 1.027106159898447
$ws.Range("K13").Value = 1.026155794721628
$ws.Range("L13").Value = 1.025136968159252
$ws.Range("M13").Value = 1.03356444541538
$ws.Range("N13").Value = 1.028564768456226

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020531801094534
$ws.Range("D14").Value = 1.022647077295699
$ws.Range("E14").Value = 1.021622095453327
$ws.Range("F14").Value = 1.030146844981729
$ws.Range("I14").Value = 1.031371721797058
$ws.Range("J14").Value = 1.027213030444705
$ws.Range("K14").Value = 1.026268175322289
$ws.Range("L14").Value = 1.025247112744368
$ws.Range("M14").Value = 1.033739528750798
$ws.Range("N14").Value = 1.028671790770918

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.020635921512866
$ws.Range("D15").Value = 1.022736440428553
$ws.Range("E15").Value = 1.021710101658498
$ws.Range("F15").Value = 1.030274805563025
$ws.Range("I15").Value = 1.031392918321396
$ws.Range("J15").Value = 1.027278861336471
$ws.Range("K15").Value = 1.02633740731767
$ws.Range("L15").Value = 1.025314969370358
$ws.Range("M15").Value = 1.033847385006211
$ws.Range("N15").Value = 1.028737715150101

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021241921556003
$ws.Range("D16").Value = 1.023256669361154
$ws.Range("E16").Value = 1.022222469255188
$ws.Range("F16").Value = 1.031019517574213
$ws.Range("I16").Value = 1.031515536502523
$ws.Range("J16").Value = 1.027661787420843
$ws.Range("K16").Value = 1.026740223558277
$ws.Range("L16").Value = 1.025709815672296
$ws.Range("M16").Value = 1.034474875048064
$ws.Range("N16").Value = 1.029121185033435

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021622026584341
$ws.Range("D17").Value = 1.023583080005202
$ws.Range("E17").Value = 1.022543980857335
$ws.Range("F17").Value = 1.031486591408643
$ws.Range("I17").Value = 1.0315917866293
$ws.Range("J17").Value = 1.027901776774571
$ws.Range("K17").Value = 1.026992772442778
$ws.Range("L17").Value = 1.025957396298222
$ws.Range("M17").Value = 1.03486823687221
$ws.Range("N17").Value = 1.029361515199572

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.021843727173711
$ws.Range("D18").Value = 1.023773500147332
$ws.Range("E18").Value = 1.022731555046827
$ws.Range("F18").Value = 1.031759004774399
$ws.Range("I18").Value = 1.031636022268086
$ws.Range("J18").Value = 1.028041682679405
$ws.Range("K18").Value = 1.027140033983202
$ws.Range("L18").Value = 1.026101771127185
$ws.Range("M18").Value = 1.035097589669509
$ws.Range("N18").Value = 1.029501619786839

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.021919319944405
$ws.Range("D19").Value = 1.023838433619764
$ws.Range("E19").Value = 1.022795520084801
$ws.Range("F19").Value = 1.031851886835148
$ws.Range("I19").Value = 1.031651064798416
$ws.Range("J19").Value = 1.028089374104685
$ws.Range("K19").Value = 1.027190238554945
$ws.Range("L19").Value = 1.026150993302068
$ws.Range("M19").Value = 1.0351757780534
$ws.Range("N19").Value = 1.029549378939414

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021581245763864
$ws.Range("D20").Value = 1.023548056089469
$ws.Range("E20").Value = 1.022509481364744
$ws.Range("F20").Value = 1.031436481181361
$ws.Range("I20").Value = 1.03158363051041
$ws.Range("J20").Value = 1.027876036028264
$ws.Range("K20").Value = 1.026965681087155
$ws.Range("L20").Value = 1.02593083683213
$ws.Range("M20").Value = 1.034826042029252
$ws.Range("N20").Value = 1.029335737898453

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020482037248431
$ws.Range("D21").Value = 1.022604368730545
$ws.Range("E21").Value = 1.021580036065559
$ws.Range("F21").Value = 1.030085686089665
$ws.Range("I21").Value = 1.03136157779063
$ws.Range("J21").Value = 1.027181562977304
$ws.Range("K21").Value = 1.0262350839942
$ws.Range("L21").Value = 1.025214679385102
$ws.Range("M21").Value = 1.033687974899515
$ws.Range("N21").Value = 1.028640278616104

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019790847601894
$ws.Range("D22").Value = 1.022011312961308
$ws.Range("E22").Value = 1.020996039487124
$ws.Range("F22").Value = 1.029236175155994
$ws.Range("I22").Value = 1.031219806622614
$ws.Range("J22").Value = 1.026744239686812
$ws.Range("K22").Value = 1.025775317736734
$ws.Range("L22").Value = 1.024764092900781
$ws.Range("M22").Value = 1.032971625506369
$ws.Range("N22").Value = 1.028202334276378

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020157267260331
$ws.Range("D23").Value = 1.022325676477993
$ws.Range("E23").Value = 1.021305590619042
$ws.Range("F23").Value = 1.029686537412428
$ws.Range("I23").Value = 1.031295166627805
$ws.Range("J23").Value = 1.026976137418964
$ws.Range("K23").Value = 1.026019087325186
$ws.Range("L23").Value = 1.025002986535562
$ws.Range("M23").Value = 1.033351452040716
$ws.Range("N23").Value = 1.028434561329908

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021599672886585
$ws.Range("D24").Value = 1.023563881795659
$ws.Range("E24").Value = 1.022525070070937
$ws.Range("F24").Value = 1.031459123906487
$ws.Range("I24").Value = 1.031587316650283
$ws.Range("J24").Value = 1.027887667397758
$ws.Range("K24").Value = 1.026977922646572
$ws.Range("L24").Value = 1.025942838018998
$ws.Range("M24").Value = 1.034845108335673
$ws.Range("N24").Value = 1.029347385785826

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023273180495377
$ws.Range("D25").Value = 1.025001937594149
$ws.Range("E25").Value = 1.023941849451731
$ws.Range("F25").Value = 1.033515216543537
$ws.Range("I25").Value = 1.031916933480201
$ws.Range("J25").Value = 1.028942471252853
$ws.Range("K25").Value = 1.028088793338981
$ws.Range("L25").Value = 1.027032117882674
$ws.Range("M25").Value = 1.036574943606848
$ws.Range("N25").Value = 1.030403687583383
